# Fix X_sqrt bug to X
# For each of the three "theta-0.2" sheets, the last existing data row (row 13)
# needs to be duplicated into three additional rows (14, 15, 16), extending the
# sheetData / dimension from E1:M13 to E1:M16.

$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "axis-1,0,0 theta-0.2",
    "axis-0,1,0 theta-0.2",
    "axis-0,0,1 theta-0.2"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $source = $ws.Range("E13:M13")

    $source.Copy($ws.Range("E14:M14"))
    $source.Copy($ws.Range("E15:M15"))
    $source.Copy($ws.Range("E16:M16"))
}
